$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H, matching the formatting of the
# existing header row (copy G1's format into H1, then set the new text)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the data value for the new Save column
$ws.Range("H2").Value = 1
